# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts and a corrected date range
# (E column) across the four worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 152
$ws.Range("F6").Value = 1076
$ws.Range("F7").Value = 2217
$ws.Range("F8").Value = 2129
$ws.Range("F9").Value = 1113
$ws.Range("F12").Value = 1676
$ws.Range("F13").Value = 400
$ws.Range("F17").Value = 216
$ws.Range("F18").Value = 1589
$ws.Range("E19").Value = "2024.06.01 10:00-06.01 18:00"
$ws.Range("F19").Value = 636
$ws.Range("F21").Value = 609
$ws.Range("F22").Value = 12272
$ws.Range("F23").Value = 12330
$ws.Range("F25").Value = 703
$ws.Range("F28").Value = 26
$ws.Range("F29").Value = 373
$ws.Range("F34").Value = 591

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 10
$ws.Range("F6").Value = 41

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 75

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 75
$ws.Range("F4").Value = 152
$ws.Range("F7").Value = 1076
$ws.Range("F8").Value = 2217
$ws.Range("F9").Value = 2129
$ws.Range("F10").Value = 1113
$ws.Range("F13").Value = 1676
$ws.Range("F14").Value = 400
$ws.Range("F18").Value = 10
$ws.Range("F21").Value = 216
$ws.Range("F22").Value = 1589
$ws.Range("E23").Value = "2024.06.01 10:00-06.01 18:00"
$ws.Range("F23").Value = 636
$ws.Range("F25").Value = 609
$ws.Range("F26").Value = 12272
$ws.Range("F27").Value = 12330
$ws.Range("F29").Value = 703
$ws.Range("F32").Value = 26
$ws.Range("F33").Value = 373
$ws.Range("F38").Value = 41
$ws.Range("F40").Value = 591

$wb.Save()
